$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 197
$ws1.Range("F5").Value = 5201
$ws1.Range("F6").Value = 28
$ws1.Range("F13").Value = 1436
$ws1.Range("F14").Value = 4142
$ws1.Range("F15").Value = 427
$ws1.Range("F16").Value = 168
$ws1.Range("F19").Value = 3167
$ws1.Range("F21").Value = 1060
$ws1.Range("F25").Value = 96
$ws1.Range("F30").Value = 19
$ws1.Range("F31").Value = 50
$ws1.Range("F34").Value = 8

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 197
$ws4.Range("F6").Value = 5201
$ws4.Range("F7").Value = 28
$ws4.Range("F14").Value = 1436
$ws4.Range("F15").Value = 4142
$ws4.Range("F16").Value = 427
$ws4.Range("F17").Value = 168
$ws4.Range("F20").Value = 3167
$ws4.Range("F22").Value = 1060
$ws4.Range("F26").Value = 96
$ws4.Range("F31").Value = 19
$ws4.Range("F32").Value = 50
$ws4.Range("F35").Value = 8
